# Update the "K" column (column G) values for rows 2-36 on Sheet1.
# The column was regenerated to report strikeouts (K) instead of the
# previous "Strike#" pitch-count metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 6
    4  = 6
    5  = 4
    6  = 6
    7  = 8
    8  = 9
    9  = 8
    10 = 5
    11 = 7
    12 = 10
    13 = 5
    14 = 2
    15 = 5
    16 = 10
    17 = 6
    18 = 8
    19 = 6
    20 = 9
    21 = 12
    22 = 6
    23 = 8
    24 = 8
    25 = 2
    26 = 2
    27 = 4
    28 = 11
    29 = 6
    30 = 5
    31 = 11
    32 = 7
    33 = 6
    34 = 9
    35 = 3
    36 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
